# The underlying change (per the diff) swaps the field data between row 7 and
# row 8 for a specific subset of columns only:
#   A  (Id), B (Taxonsorteringsordning), E (TaxonId), F (Artnamn),
#   G  (Vetenskapligt namn), H (Auktor), I (Antal), Q (Ost), R (Nord),
#   AC (Publik kommentar)
# All other columns (C, D, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT,
# AW, AX, AY, ...) stay exactly as they were on both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column "I" ("Antal") holds values that look like plain numbers (e.g. "20")
# but are actually stored as text in the source file. A normal
# Range.Value2 assignment of a numeric-looking string gets auto-coerced to a
# real number by Excel, which would change the stored cell type. To preserve
# the original text typing we "launder" the swap through Copy +
# PasteSpecial(xlPasteValues), which transfers the cell's value (including
# its text-vs-number flavor) without touching formatting/styles.
$ws.Range("I7").Copy()
$ws.Range("I8").PasteSpecial(-4163)  # -4163 == xlPasteValues
$ws.Range("I7").ClearContents()

# Remaining swapped columns are plain numbers or free-text strings that do
# not round-trip ambiguously, so a direct Value2 swap is safe.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

foreach ($col in $cols) {
    $addr7 = "$col" + "7"
    $addr8 = "$col" + "8"

    $range7 = $ws.Range($addr7)
    $range8 = $ws.Range($addr8)

    $val7 = $range7.Value2
    $val8 = $range8.Value2

    $range7.Value2 = $val8
    $range8.Value2 = $val7
}
